# Add "rebounds", "assists" and "personalfouls" columns (L, M, N) to the
# Lakers sheet header row, mirroring the existing H/I/J headers (these new
# columns feed the already-present L:N ROUNDUP(...)*10 helper formulas that
# power the new "functioning lakers api route").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lakers")

$ws.Range("L1").Value = "rebounds"
$ws.Range("M1").Value = "assists"
$ws.Range("N1").Value = "personalfouls"

# Column widths adjust as a side effect of the new header text needing to
# fit (best-fit sizing). Reproduce the resulting widths as closely as this
# runtime's column-width model allows.
$ws.Columns.Item(7).ColumnWidth = 5.666666666666667
$ws.Columns.Item(10).ColumnWidth = 12.5
$ws.Columns.Item(11).ColumnWidth = 17
$ws.Columns.Item(13).ColumnWidth = 5.833333333333333
$ws.Columns.Item(14).ColumnWidth = 12.5

# The active selection moved one column to the right (R13 -> S13).
$ws.Range("S13").Select()
